# Commit: "Muuta Tilaajan turvallisuusvastaavan oikeudet urakkakohtaisiksi."
# ("Change the Tilaajan turvallisuusvastaava's [Client's safety officer]
# permissions to be contract-specific.")
#
# On the "Oikeudet" sheet, column I holds the access rights of the role
# "Tilaajan turvallisuusvastaava". Every "R*" (read access to ALL contracts)
# becomes plain "R" (read access to named/assigned contracts only), and the
# two "R*,W*" cells become "R,W" - i.e. the trailing "*" (which the sheet's
# legend defines as "all contracts") is dropped throughout the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column I (Tilaajan turvallisuusvastaava): drop the "*" (all-contracts)
# marker, making the role's access rights contract-specific instead of
# blanket access to all contracts (R* -> R, R*,W* -> R,W). ---

$ws.Range("I7").Value = "R"
$ws.Range("I13").Value = "R"
$ws.Range("I14").Value = "R"
$ws.Range("I18").Value = "R"
$ws.Range("I19").Value = "R"
$ws.Range("I22").Value = "R"
$ws.Range("I23").Value = "R"
$ws.Range("I24").Value = "R"
$ws.Range("I25").Value = "R"
$ws.Range("I26").Value = "R"
$ws.Range("I28").Value = "R"
$ws.Range("I50").Value = "R"
$ws.Range("I56").Value = "R"
$ws.Range("I57").Value = "R"
$ws.Range("I63").Value = "R"
$ws.Range("I64").Value = "R"
$ws.Range("I65").Value = "R"
$ws.Range("I66").Value = "R"
$ws.Range("I67").Value = "R"
$ws.Range("I68").Value = "R"
$ws.Range("I69").Value = "R"
$ws.Range("I70").Value = "R"
$ws.Range("I71").Value = "R"
$ws.Range("I72").Value = "R"
$ws.Range("I77").Value = "R"
$ws.Range("I78").Value = "R"
$ws.Range("I79").Value = "R"
$ws.Range("I80").Value = "R"
$ws.Range("I81").Value = "R"
$ws.Range("I82").Value = "R"
$ws.Range("I83").Value = "R"
$ws.Range("I84").Value = "R"
$ws.Range("I86").Value = "R"
$ws.Range("I89").Value = "R"
$ws.Range("I90").Value = "R"
$ws.Range("I91").Value = "R"
$ws.Range("I92").Value = "R"
$ws.Range("I93").Value = "R"

$ws.Range("I29").Value = "R,W"
$ws.Range("I45").Value = "R,W"

# These rows also pick up a thin left border when Excel normalizes the
# cell style on edit (matching the neighbouring, already-bordered cells).
$ws.Range("I21").Value = "R"
$ws.Range("I21").Borders.Item(7).LineStyle = 1
$ws.Range("I21").Borders.Item(7).Weight = 2
$ws.Range("I27").Value = "R"
$ws.Range("I27").Borders.Item(7).LineStyle = 1
$ws.Range("I27").Borders.Item(7).Weight = 2
$ws.Range("I36").Value = "R"
$ws.Range("I36").Borders.Item(7).LineStyle = 1
$ws.Range("I36").Borders.Item(7).Weight = 2
$ws.Range("I37").Value = "R"
$ws.Range("I37").Borders.Item(7).LineStyle = 1
$ws.Range("I37").Borders.Item(7).Weight = 2
$ws.Range("I38").Value = "R"
$ws.Range("I38").Borders.Item(7).LineStyle = 1
$ws.Range("I38").Borders.Item(7).Weight = 2
$ws.Range("I39").Value = "R"
$ws.Range("I39").Borders.Item(7).LineStyle = 1
$ws.Range("I39").Borders.Item(7).Weight = 2
$ws.Range("I40").Value = "R"
$ws.Range("I40").Borders.Item(7).LineStyle = 1
$ws.Range("I40").Borders.Item(7).Weight = 2
$ws.Range("I41").Value = "R"
$ws.Range("I41").Borders.Item(7).LineStyle = 1
$ws.Range("I41").Borders.Item(7).Weight = 2
$ws.Range("I42").Value = "R"
$ws.Range("I42").Borders.Item(7).LineStyle = 1
$ws.Range("I42").Borders.Item(7).Weight = 2
$ws.Range("I43").Value = "R"
$ws.Range("I43").Borders.Item(7).LineStyle = 1
$ws.Range("I43").Borders.Item(7).Weight = 2
$ws.Range("I44").Value = "R"
$ws.Range("I44").Borders.Item(7).LineStyle = 1
$ws.Range("I44").Borders.Item(7).Weight = 2

# The column is now much longer ("R*,W*,aseta-näkyviin-urakoitsijalle" etc.
# still live further down as other roles), so the author widened columns I
# and P (which holds the analogous "Tilaajan_Kayttaja" rights) and let Excel
# drop their old best-fit auto-size flag.
$ws.Columns.Item(9).ColumnWidth = 36
$ws.Columns.Item(16).ColumnWidth = 29

# Restore/record the cursor position on both sheets (sheet2's selection is
# set first so that sheet1 - the tab that was active before the edit -
# ends up the active tab again).
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G36").Select()

$ws.Range("I82").Select()
